$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 390 (shifts existing rows 390..446 down to 391..447,
# formatting is inherited from the row above - matches the style used by column D
# in this data block).
$ws.Rows.Item(390).Insert()

# Populate the newly inserted row 390 with the new weekly data record.
$ws.Cells.Item(390, 1).Value  = 8
$ws.Cells.Item(390, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(390, 3).Value  = "Coquimbo"
$ws.Cells.Item(390, 4).Value  = 45127
$ws.Cells.Item(390, 5).Value  = 4
$ws.Cells.Item(390, 6).Value  = 100112012
$ws.Cells.Item(390, 7).Value  = "Espinaca"
$ws.Cells.Item(390, 8).Value  = "Sin especificar"
$ws.Cells.Item(390, 9).Value  = "Primera"
$ws.Cells.Item(390, 10).Value = 1600
$ws.Cells.Item(390, 11).Value = 500
$ws.Cells.Item(390, 12).Value = 600
$ws.Cells.Item(390, 13).Value = 550
$ws.Cells.Item(390, 14).Value = "$/atado 300 a 500 gramos"
$ws.Cells.Item(390, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(390, 16).Value = 1100
$ws.Cells.Item(390, 17).Value = 0.5
$ws.Cells.Item(390, 18).Value = "Hortaliza"
